$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LIST")

# Fill in the role/profile codes in column A (rows 2-8), sourced from the
# reference list on Feuil1 (B9:B14) plus "MP.CPT" inserted before the last
# entry. This also removes the stray "RO.FOU.001.CRE.01" value that used
# to sit in A3.
$ws.Range("A2").Value = "AD.SEC.001.FON.02"
$ws.Range("A3").Value = "AD.SEC.001.FON.01"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A4").Value = "AD.SEC.001.FON.03"
$ws.Range("A5").Value = "RO.ACT"
$ws.Range("A6").Value = "RO.FOU"
$ws.Range("A7").Value = "MP.CPT"
$ws.Range("A8").Value = "AD.SEC.014.FON.01"

# Move the active selection to B15, matching the saved view state.
$ws.Activate()
$ws.Range("B15").Select()
